$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Revert "new changes in ops (ordercreation & orderpage & order form)":
# the sheet is collapsed back from 15 columns (A:O) to 13 columns (A:M) by
# removing the "Typist" / "Typist QC" columns, and the Search(T1) tier
# value replaces Typing(T1).
# ---------------------------------------------------------------------------

# --- fix up formatting first (copy formats before the values are overwritten) ---
# The highlighted "State/County" style (style index 4) moves from K:L to I:J.
$ws.Range("K2:L3").Copy()
$ws.Range("I2:J3").PasteSpecial(-4122)

# K:L become plain data columns again (style index 1, like any other text cell).
$ws.Range("B2:B3").Copy()
$ws.Range("K2:L3").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- header row (row 1) ---
$ws.Range("E1").Value = "Client"
$ws.Range("F1").Value = "Lob"
$ws.Range("G1").Value = "Process"
$ws.Range("H1").Value = "Product Name"
$ws.Range("I1").Value = "State"
$ws.Range("J1").Value = "County"
$ws.Range("K1").Value = "Municipality"
$ws.Range("L1").Value = "Status"
$ws.Range("M1").Value = "Tier"

# --- row 2 ---
$ws.Range("E2").Value = "Flowers Title Companies"
$ws.Range("F2").Value = "Title"
$ws.Range("G2").Value = "Production & QC"
$ws.Range("H2").Value = "Order Entry – ETTC"
$ws.Range("I2").Value = "AL"
$ws.Range("J2").Value = "Autauga"
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = "WIP"
$ws.Range("M2").Value = "Search(T1)"

# --- row 3 ---
$ws.Range("C3").Value = "SIPL5317"
$ws.Range("D3").Value = "SIPL5688"
$ws.Range("E3").Value = "Flowers Title Companies"
$ws.Range("F3").Value = "Title"
$ws.Range("G3").Value = "Typing"
$ws.Range("H3").Value = "Policy Typing"
$ws.Range("I3").Value = "AL"
$ws.Range("J3").Value = "Baldwin"
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = "WIP"
$ws.Range("M3").Value = "Search(T1)"

# the old N:O columns (Status/Tier) no longer exist, clear their leftovers
$ws.Range("N1:O3").Clear()

# --- re-fit the columns whose best-fit width changed because of the new text ---
$ws.Columns("E").ColumnWidth = 19.833333333333336
$ws.Columns("H").ColumnWidth = 33.666666666666664
$ws.Columns("J").ColumnWidth = 12

# --- selection state matches the reverted file ---
$ws.Range("A4:XFD5").Select()

# ---------------------------------------------------------------------------
# Re-create the (now orphaned) conditional-formatting differential style that
# ships with the reverted workbook: add then remove a "duplicate values"
# rule so the dxf record is written to styles.xml without leaving a live
# conditionalFormatting block on the sheet.
# ---------------------------------------------------------------------------
$dxfRange = $ws.Range("A1:A1")
$fc = $dxfRange.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615
$dxfRange.FormatConditions.Delete()

Write-Host "edit applied"
